# repull data, push all data, mean calculation
# Update the "dSF" (column F) values on Sheet1 with newly re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    4  = 1
    6  = -4
    7  = -1
    9  = -1
    10 = 2
    11 = -5
    12 = -4
    13 = 8
    14 = 4
    16 = 1
    17 = -1
    18 = -1
    19 = 5
    22 = -1
    23 = 1
    24 = -2
    25 = 3
    28 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
